# Spon2-Itga4.xlsx : refresh LR-pairs TPM-derived metrics (cols E:T) for all 30 data rows.
# (Sending cluster / Ligand / Receptor / Target cluster labels in A:D are unchanged.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 30  # sheet rows 2..31
$cols = 20  # columns A..T
$arr = New-Object 'object[,]' $rows,$cols

# row 2: ECs -> ECs
$arr[0,0] = "ECs"
$arr[0,1] = "Spon2"
$arr[0,2] = "Itga4"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = [double]"0.7021656666666667"
$arr[0,7] = [double]"2.106497"
$arr[0,8] = [double]"0.197084787652417"
$arr[0,9] = [double]"0.197084787652417"
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = [double]"1.105124"
$arr[0,13] = [double]"3.315372"
$arr[0,14] = [double]"0.006910839970832482"
$arr[0,15] = [double]"0.006910839970832482"
$arr[0,16] = [double]"0.7759801302093333"
$arr[0,17] = [double]"6.983821171884"
$arr[0,18] = [double]"0.001362021428151355"
$arr[0,19] = [double]"0.001362021428151355"

# row 3: ECs -> FAPs
$arr[1,0] = "ECs"
$arr[1,1] = "Spon2"
$arr[1,2] = "Itga4"
$arr[1,3] = "FAPs"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = [double]"0.7021656666666667"
$arr[1,7] = [double]"2.106497"
$arr[1,8] = [double]"0.197084787652417"
$arr[1,9] = [double]"0.197084787652417"
$arr[1,10] = 1
$arr[1,11] = [double]"0.3333333333333333"
$arr[1,12] = [double]"0.04442266666666667"
$arr[1,13] = [double]"0.133268"
$arr[1,14] = [double]"0.0002777950170396876"
$arr[1,15] = [double]"0.0002777950170396876"
$arr[1,16] = [double]"0.03119207135511111"
$arr[1,17] = [double]"0.280728642196"
$arr[1,18] = [double]"5.474917194416638e-05"
$arr[1,19] = [double]"5.474917194416638e-05"

# row 4: ECs -> Inflammatory-Mac
$arr[2,0] = "ECs"
$arr[2,1] = "Spon2"
$arr[2,2] = "Itga4"
$arr[2,3] = "Inflammatory-Mac"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = [double]"0.7021656666666667"
$arr[2,7] = [double]"2.106497"
$arr[2,8] = [double]"0.197084787652417"
$arr[2,9] = [double]"0.197084787652417"
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = [double]"56.54517366666666"
$arr[2,13] = [double]"169.635521"
$arr[2,14] = [double]"0.3536025335919447"
$arr[2,15] = [double]"0.3536025335919447"
$arr[2,16] = [double]"39.70407956443744"
$arr[2,17] = [double]"357.336716079937"
$arr[2,18] = [double]"0.06968968024632505"
$arr[2,19] = [double]"0.06968968024632505"

# row 5: ECs -> MuSCs
$arr[3,0] = "ECs"
$arr[3,1] = "Spon2"
$arr[3,2] = "Itga4"
$arr[3,3] = "MuSCs"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = [double]"0.7021656666666667"
$arr[3,7] = [double]"2.106497"
$arr[3,8] = [double]"0.197084787652417"
$arr[3,9] = [double]"0.197084787652417"
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = [double]"0.8044289999999998"
$arr[3,13] = [double]"2.413287"
$arr[3,14] = [double]"0.005030458199167516"
$arr[3,15] = [double]"0.005030458199167516"
$arr[3,16] = [double]"0.5648424250709999"
$arr[3,17] = [double]"5.083581825639"
$arr[3,18] = [double]"0.0009914267859772897"
$arr[3,19] = [double]"0.0009914267859772897"

# row 6: ECs -> Neutrophils
$arr[4,0] = "ECs"
$arr[4,1] = "Spon2"
$arr[4,2] = "Itga4"
$arr[4,3] = "Neutrophils"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = [double]"0.7021656666666667"
$arr[4,7] = [double]"2.106497"
$arr[4,8] = [double]"0.197084787652417"
$arr[4,9] = [double]"0.197084787652417"
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = [double]"78.08909333333334"
$arr[4,13] = [double]"234.26728"
$arr[4,14] = [double]"0.4883264027331488"
$arr[4,15] = [double]"0.4883264027331488"
$arr[4,16] = [double]"54.83148027979556"
$arr[4,17] = [double]"493.48332251816"
$arr[4,18] = [double]"0.09624170538773127"
$arr[4,19] = [double]"0.09624170538773127"

# row 7: ECs -> Resolving-Mac
$arr[5,0] = "ECs"
$arr[5,1] = "Spon2"
$arr[5,2] = "Itga4"
$arr[5,3] = "Resolving-Mac"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = [double]"0.7021656666666667"
$arr[5,7] = [double]"2.106497"
$arr[5,8] = [double]"0.197084787652417"
$arr[5,9] = [double]"0.197084787652417"
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = [double]"23.323433"
$arr[5,13] = [double]"69.970299"
$arr[5,14] = [double]"0.1458519704878668"
$arr[5,15] = [double]"0.1458519704878668"
$arr[5,16] = [double]"16.37691388140033"
$arr[5,17] = [double]"147.392224932603"
$arr[5,18] = [double]"0.02874520463228782"
$arr[5,19] = [double]"0.02874520463228782"

# row 8: FAPs -> ECs
$arr[6,0] = "FAPs"
$arr[6,1] = "Spon2"
$arr[6,2] = "Itga4"
$arr[6,3] = "ECs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = [double]"2.541923333333334"
$arr[6,7] = [double]"7.62577"
$arr[6,8] = [double]"0.7134704018739033"
$arr[6,9] = [double]"0.7134704018739033"
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = [double]"1.105124"
$arr[6,13] = [double]"3.315372"
$arr[6,14] = [double]"0.006910839970832482"
$arr[6,15] = [double]"0.006910839970832482"
$arr[6,16] = [double]"2.809140481826667"
$arr[6,17] = [double]"25.28226433644"
$arr[6,18] = [double]"0.004930679771276086"
$arr[6,19] = [double]"0.004930679771276086"

# row 9: FAPs -> FAPs
$arr[7,0] = "FAPs"
$arr[7,1] = "Spon2"
$arr[7,2] = "Itga4"
$arr[7,3] = "FAPs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = [double]"2.541923333333334"
$arr[7,7] = [double]"7.62577"
$arr[7,8] = [double]"0.7134704018739033"
$arr[7,9] = [double]"0.7134704018739033"
$arr[7,10] = 1
$arr[7,11] = [double]"0.3333333333333333"
$arr[7,12] = [double]"0.04442266666666667"
$arr[7,13] = [double]"0.133268"
$arr[7,14] = [double]"0.0002777950170396876"
$arr[7,15] = [double]"0.0002777950170396876"
$arr[7,16] = [double]"0.1129190129288889"
$arr[7,17] = [double]"1.01627111636"
$arr[7,18] = [double]"0.0001981985224458737"
$arr[7,19] = [double]"0.0001981985224458737"

# row 10: FAPs -> Inflammatory-Mac
$arr[8,0] = "FAPs"
$arr[8,1] = "Spon2"
$arr[8,2] = "Itga4"
$arr[8,3] = "Inflammatory-Mac"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = [double]"2.541923333333334"
$arr[8,7] = [double]"7.62577"
$arr[8,8] = [double]"0.7134704018739033"
$arr[8,9] = [double]"0.7134704018739033"
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = [double]"56.54517366666666"
$arr[8,13] = [double]"169.635521"
$arr[8,14] = [double]"0.3536025335919447"
$arr[8,15] = [double]"0.3536025335919447"
$arr[8,16] = [double]"143.7334963306856"
$arr[8,17] = [double]"1293.60146697617"
$arr[8,18] = [double]"0.2522849417454752"
$arr[8,19] = [double]"0.2522849417454752"

# row 11: FAPs -> MuSCs
$arr[9,0] = "FAPs"
$arr[9,1] = "Spon2"
$arr[9,2] = "Itga4"
$arr[9,3] = "MuSCs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = [double]"2.541923333333334"
$arr[9,7] = [double]"7.62577"
$arr[9,8] = [double]"0.7134704018739033"
$arr[9,9] = [double]"0.7134704018739033"
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = [double]"0.8044289999999998"
$arr[9,13] = [double]"2.413287"
$arr[9,14] = [double]"0.005030458199167516"
$arr[9,15] = [double]"0.005030458199167516"
$arr[9,16] = [double]"2.04479684511"
$arr[9,17] = [double]"18.40317160599"
$arr[9,18] = [double]"0.003589083032969919"
$arr[9,19] = [double]"0.003589083032969919"

# row 12: FAPs -> Neutrophils
$arr[10,0] = "FAPs"
$arr[10,1] = "Spon2"
$arr[10,2] = "Itga4"
$arr[10,3] = "Neutrophils"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = [double]"2.541923333333334"
$arr[10,7] = [double]"7.62577"
$arr[10,8] = [double]"0.7134704018739033"
$arr[10,9] = [double]"0.7134704018739033"
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = [double]"78.08909333333334"
$arr[10,13] = [double]"234.26728"
$arr[10,14] = [double]"0.4883264027331488"
$arr[10,15] = [double]"0.4883264027331488"
$arr[10,16] = [double]"198.4964884228445"
$arr[10,17] = [double]"1786.4683958056"
$arr[10,18] = [double]"0.3484064348036572"
$arr[10,19] = [double]"0.3484064348036572"

# row 13: FAPs -> Resolving-Mac
$arr[11,0] = "FAPs"
$arr[11,1] = "Spon2"
$arr[11,2] = "Itga4"
$arr[11,3] = "Resolving-Mac"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = [double]"2.541923333333334"
$arr[11,7] = [double]"7.62577"
$arr[11,8] = [double]"0.7134704018739033"
$arr[11,9] = [double]"0.7134704018739033"
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = [double]"23.323433"
$arr[11,13] = [double]"69.970299"
$arr[11,14] = [double]"0.1458519704878668"
$arr[11,15] = [double]"0.1458519704878668"
$arr[11,16] = [double]"59.28637855613667"
$arr[11,17] = [double]"533.57740700523"
$arr[11,18] = [double]"0.104061063998079"
$arr[11,19] = [double]"0.104061063998079"

# row 14: Inflammatory-Mac -> ECs
$arr[12,0] = "Inflammatory-Mac"
$arr[12,1] = "Spon2"
$arr[12,2] = "Itga4"
$arr[12,3] = "ECs"
$arr[12,4] = 2
$arr[12,5] = [double]"0.6666666666666666"
$arr[12,6] = [double]"0.1068446666666667"
$arr[12,7] = [double]"0.320534"
$arr[12,8] = [double]"0.0299893022992104"
$arr[12,9] = [double]"0.0299893022992104"
$arr[12,10] = 3
$arr[12,11] = 1
$arr[12,12] = [double]"1.105124"
$arr[12,13] = [double]"3.315372"
$arr[12,14] = [double]"0.006910839970832482"
$arr[12,15] = [double]"0.006910839970832482"
$arr[12,16] = [double]"0.1180766054053333"
$arr[12,17] = [double]"1.062689448648"
$arr[12,18] = [double]"0.0002072512690267617"
$arr[12,19] = [double]"0.0002072512690267617"

# row 15: Inflammatory-Mac -> FAPs
$arr[13,0] = "Inflammatory-Mac"
$arr[13,1] = "Spon2"
$arr[13,2] = "Itga4"
$arr[13,3] = "FAPs"
$arr[13,4] = 2
$arr[13,5] = [double]"0.6666666666666666"
$arr[13,6] = [double]"0.1068446666666667"
$arr[13,7] = [double]"0.320534"
$arr[13,8] = [double]"0.0299893022992104"
$arr[13,9] = [double]"0.0299893022992104"
$arr[13,10] = 1
$arr[13,11] = [double]"0.3333333333333333"
$arr[13,12] = [double]"0.04442266666666667"
$arr[13,13] = [double]"0.133268"
$arr[13,14] = [double]"0.0002777950170396876"
$arr[13,15] = [double]"0.0002777950170396876"
$arr[13,16] = [double]"0.004746325012444444"
$arr[13,17] = [double]"0.042716925112"
$arr[13,18] = [double]"8.330878743217496e-06"
$arr[13,19] = [double]"8.330878743217496e-06"

# row 16: Inflammatory-Mac -> Inflammatory-Mac
$arr[14,0] = "Inflammatory-Mac"
$arr[14,1] = "Spon2"
$arr[14,2] = "Itga4"
$arr[14,3] = "Inflammatory-Mac"
$arr[14,4] = 2
$arr[14,5] = [double]"0.6666666666666666"
$arr[14,6] = [double]"0.1068446666666667"
$arr[14,7] = [double]"0.320534"
$arr[14,8] = [double]"0.0299893022992104"
$arr[14,9] = [double]"0.0299893022992104"
$arr[14,10] = 3
$arr[14,11] = 1
$arr[14,12] = [double]"56.54517366666666"
$arr[14,13] = [double]"169.635521"
$arr[14,14] = [double]"0.3536025335919447"
$arr[14,15] = [double]"0.3536025335919447"
$arr[14,16] = [double]"6.041550232023777"
$arr[14,17] = [double]"54.37395208821399"
$arr[14,18] = [double]"0.01060429327365553"
$arr[14,19] = [double]"0.01060429327365553"

# row 17: Inflammatory-Mac -> MuSCs
$arr[15,0] = "Inflammatory-Mac"
$arr[15,1] = "Spon2"
$arr[15,2] = "Itga4"
$arr[15,3] = "MuSCs"
$arr[15,4] = 2
$arr[15,5] = [double]"0.6666666666666666"
$arr[15,6] = [double]"0.1068446666666667"
$arr[15,7] = [double]"0.320534"
$arr[15,8] = [double]"0.0299893022992104"
$arr[15,9] = [double]"0.0299893022992104"
$arr[15,10] = 3
$arr[15,11] = 1
$arr[15,12] = [double]"0.8044289999999998"
$arr[15,13] = [double]"2.413287"
$arr[15,14] = [double]"0.005030458199167516"
$arr[15,15] = [double]"0.005030458199167516"
$arr[15,16] = [double]"0.08594894836199997"
$arr[15,17] = [double]"0.7735405352579998"
$arr[15,18] = [double]"0.0001508599316383762"
$arr[15,19] = [double]"0.0001508599316383762"

# row 18: Inflammatory-Mac -> Neutrophils
$arr[16,0] = "Inflammatory-Mac"
$arr[16,1] = "Spon2"
$arr[16,2] = "Itga4"
$arr[16,3] = "Neutrophils"
$arr[16,4] = 2
$arr[16,5] = [double]"0.6666666666666666"
$arr[16,6] = [double]"0.1068446666666667"
$arr[16,7] = [double]"0.320534"
$arr[16,8] = [double]"0.0299893022992104"
$arr[16,9] = [double]"0.0299893022992104"
$arr[16,10] = 3
$arr[16,11] = 1
$arr[16,12] = [double]"78.08909333333334"
$arr[16,13] = [double]"234.26728"
$arr[16,14] = [double]"0.4883264027331488"
$arr[16,15] = [double]"0.4883264027331488"
$arr[16,16] = [double]"8.343403147502222"
$arr[16,17] = [double]"75.09062832752"
$arr[16,18] = [double]"0.01464456811225036"
$arr[16,19] = [double]"0.01464456811225036"

# row 19: Inflammatory-Mac -> Resolving-Mac
$arr[17,0] = "Inflammatory-Mac"
$arr[17,1] = "Spon2"
$arr[17,2] = "Itga4"
$arr[17,3] = "Resolving-Mac"
$arr[17,4] = 2
$arr[17,5] = [double]"0.6666666666666666"
$arr[17,6] = [double]"0.1068446666666667"
$arr[17,7] = [double]"0.320534"
$arr[17,8] = [double]"0.0299893022992104"
$arr[17,9] = [double]"0.0299893022992104"
$arr[17,10] = 3
$arr[17,11] = 1
$arr[17,12] = [double]"23.323433"
$arr[17,13] = [double]"69.970299"
$arr[17,14] = [double]"0.1458519704878668"
$arr[17,15] = [double]"0.1458519704878668"
$arr[17,16] = [double]"2.491984424407333"
$arr[17,17] = [double]"22.427859819666"
$arr[17,18] = [double]"0.004373998833896152"
$arr[17,19] = [double]"0.004373998833896152"

# row 20: MuSCs -> ECs
$arr[18,0] = "MuSCs"
$arr[18,1] = "Spon2"
$arr[18,2] = "Itga4"
$arr[18,3] = "ECs"
$arr[18,4] = 1
$arr[18,5] = [double]"0.3333333333333333"
$arr[18,6] = [double]"0.08902066666666668"
$arr[18,7] = [double]"0.267062"
$arr[18,8] = [double]"0.02498643841412059"
$arr[18,9] = [double]"0.02498643841412059"
$arr[18,10] = 3
$arr[18,11] = 1
$arr[18,12] = [double]"1.105124"
$arr[18,13] = [double]"3.315372"
$arr[18,14] = [double]"0.006910839970832482"
$arr[18,15] = [double]"0.006910839970832482"
$arr[18,16] = [double]"0.09837887522933335"
$arr[18,17] = [double]"0.8854098770640001"
$arr[18,18] = [double]"0.0001726772773210488"
$arr[18,19] = [double]"0.0001726772773210488"

# row 21: MuSCs -> FAPs
$arr[19,0] = "MuSCs"
$arr[19,1] = "Spon2"
$arr[19,2] = "Itga4"
$arr[19,3] = "FAPs"
$arr[19,4] = 1
$arr[19,5] = [double]"0.3333333333333333"
$arr[19,6] = [double]"0.08902066666666668"
$arr[19,7] = [double]"0.267062"
$arr[19,8] = [double]"0.02498643841412059"
$arr[19,9] = [double]"0.02498643841412059"
$arr[19,10] = 1
$arr[19,11] = [double]"0.3333333333333333"
$arr[19,12] = [double]"0.04442266666666667"
$arr[19,13] = [double]"0.133268"
$arr[19,14] = [double]"0.0002777950170396876"
$arr[19,15] = [double]"0.0002777950170396876"
$arr[19,16] = [double]"0.003954535401777778"
$arr[19,17] = [double]"0.035590818616"
$arr[19,18] = [double]"6.941108085011735e-06"
$arr[19,19] = [double]"6.941108085011735e-06"

# row 22: MuSCs -> Inflammatory-Mac
$arr[20,0] = "MuSCs"
$arr[20,1] = "Spon2"
$arr[20,2] = "Itga4"
$arr[20,3] = "Inflammatory-Mac"
$arr[20,4] = 1
$arr[20,5] = [double]"0.3333333333333333"
$arr[20,6] = [double]"0.08902066666666668"
$arr[20,7] = [double]"0.267062"
$arr[20,8] = [double]"0.02498643841412059"
$arr[20,9] = [double]"0.02498643841412059"
$arr[20,10] = 3
$arr[20,11] = 1
$arr[20,12] = [double]"56.54517366666666"
$arr[20,13] = [double]"169.635521"
$arr[20,14] = [double]"0.3536025335919447"
$arr[20,15] = [double]"0.3536025335919447"
$arr[20,16] = [double]"5.033689056589112"
$arr[20,17] = [double]"45.303201509302"
$arr[20,18] = [double]"0.008835267928672133"
$arr[20,19] = [double]"0.008835267928672133"

# row 23: MuSCs -> MuSCs
$arr[21,0] = "MuSCs"
$arr[21,1] = "Spon2"
$arr[21,2] = "Itga4"
$arr[21,3] = "MuSCs"
$arr[21,4] = 1
$arr[21,5] = [double]"0.3333333333333333"
$arr[21,6] = [double]"0.08902066666666668"
$arr[21,7] = [double]"0.267062"
$arr[21,8] = [double]"0.02498643841412059"
$arr[21,9] = [double]"0.02498643841412059"
$arr[21,10] = 3
$arr[21,11] = 1
$arr[21,12] = [double]"0.8044289999999998"
$arr[21,13] = [double]"2.413287"
$arr[21,14] = [double]"0.005030458199167516"
$arr[21,15] = [double]"0.005030458199167516"
$arr[21,16] = [double]"0.071610805866"
$arr[21,17] = [double]"0.6444972527939999"
$arr[21,18] = [double]"0.0001256932339883071"
$arr[21,19] = [double]"0.0001256932339883071"

# row 24: MuSCs -> Neutrophils
$arr[22,0] = "MuSCs"
$arr[22,1] = "Spon2"
$arr[22,2] = "Itga4"
$arr[22,3] = "Neutrophils"
$arr[22,4] = 1
$arr[22,5] = [double]"0.3333333333333333"
$arr[22,6] = [double]"0.08902066666666668"
$arr[22,7] = [double]"0.267062"
$arr[22,8] = [double]"0.02498643841412059"
$arr[22,9] = [double]"0.02498643841412059"
$arr[22,10] = 3
$arr[22,11] = 1
$arr[22,12] = [double]"78.08909333333334"
$arr[22,13] = [double]"234.26728"
$arr[22,14] = [double]"0.4883264027331488"
$arr[22,15] = [double]"0.4883264027331488"
$arr[22,16] = [double]"6.95154314792889"
$arr[22,17] = [double]"62.56388833136"
$arr[22,18] = [double]"0.01220153758788087"
$arr[22,19] = [double]"0.01220153758788087"

# row 25: MuSCs -> Resolving-Mac
$arr[23,0] = "MuSCs"
$arr[23,1] = "Spon2"
$arr[23,2] = "Itga4"
$arr[23,3] = "Resolving-Mac"
$arr[23,4] = 1
$arr[23,5] = [double]"0.3333333333333333"
$arr[23,6] = [double]"0.08902066666666668"
$arr[23,7] = [double]"0.267062"
$arr[23,8] = [double]"0.02498643841412059"
$arr[23,9] = [double]"0.02498643841412059"
$arr[23,10] = 3
$arr[23,11] = 1
$arr[23,12] = [double]"23.323433"
$arr[23,13] = [double]"69.970299"
$arr[23,14] = [double]"0.1458519704878668"
$arr[23,15] = [double]"0.1458519704878668"
$arr[23,16] = [double]"2.076267554615333"
$arr[23,17] = [double]"18.686407991538"
$arr[23,18] = [double]"0.003644321278173219"
$arr[23,19] = [double]"0.003644321278173219"

# row 26: Neutrophils -> ECs
$arr[24,0] = "Neutrophils"
$arr[24,1] = "Spon2"
$arr[24,2] = "Itga4"
$arr[24,3] = "ECs"
$arr[24,4] = 2
$arr[24,5] = [double]"0.6666666666666666"
$arr[24,6] = [double]"0.122805"
$arr[24,7] = [double]"0.368415"
$arr[24,8] = [double]"0.03446906976034867"
$arr[24,9] = [double]"0.03446906976034867"
$arr[24,10] = 3
$arr[24,11] = 1
$arr[24,12] = [double]"1.105124"
$arr[24,13] = [double]"3.315372"
$arr[24,14] = [double]"0.006910839970832482"
$arr[24,15] = [double]"0.006910839970832482"
$arr[24,16] = [double]"0.13571475282"
$arr[24,17] = [double]"1.22143277538"
$arr[24,18] = [double]"0.0002382102250572308"
$arr[24,19] = [double]"0.0002382102250572308"

# row 27: Neutrophils -> FAPs
$arr[25,0] = "Neutrophils"
$arr[25,1] = "Spon2"
$arr[25,2] = "Itga4"
$arr[25,3] = "FAPs"
$arr[25,4] = 2
$arr[25,5] = [double]"0.6666666666666666"
$arr[25,6] = [double]"0.122805"
$arr[25,7] = [double]"0.368415"
$arr[25,8] = [double]"0.03446906976034867"
$arr[25,9] = [double]"0.03446906976034867"
$arr[25,10] = 1
$arr[25,11] = [double]"0.3333333333333333"
$arr[25,12] = [double]"0.04442266666666667"
$arr[25,13] = [double]"0.133268"
$arr[25,14] = [double]"0.0002777950170396876"
$arr[25,15] = [double]"0.0002777950170396876"
$arr[25,16] = [double]"0.00545532558"
$arr[25,17] = [double]"0.04909793022"
$arr[25,18] = [double]"9.57533582141824e-06"
$arr[25,19] = [double]"9.57533582141824e-06"

# row 28: Neutrophils -> Inflammatory-Mac
$arr[26,0] = "Neutrophils"
$arr[26,1] = "Spon2"
$arr[26,2] = "Itga4"
$arr[26,3] = "Inflammatory-Mac"
$arr[26,4] = 2
$arr[26,5] = [double]"0.6666666666666666"
$arr[26,6] = [double]"0.122805"
$arr[26,7] = [double]"0.368415"
$arr[26,8] = [double]"0.03446906976034867"
$arr[26,9] = [double]"0.03446906976034867"
$arr[26,10] = 3
$arr[26,11] = 1
$arr[26,12] = [double]"56.54517366666666"
$arr[26,13] = [double]"169.635521"
$arr[26,14] = [double]"0.3536025335919447"
$arr[26,15] = [double]"0.3536025335919447"
$arr[26,16] = [double]"6.944030052134999"
$arr[26,17] = [double]"62.49627046921499"
$arr[26,18] = [double]"0.01218835039781677"
$arr[26,19] = [double]"0.01218835039781677"

# row 29: Neutrophils -> MuSCs
$arr[27,0] = "Neutrophils"
$arr[27,1] = "Spon2"
$arr[27,2] = "Itga4"
$arr[27,3] = "MuSCs"
$arr[27,4] = 2
$arr[27,5] = [double]"0.6666666666666666"
$arr[27,6] = [double]"0.122805"
$arr[27,7] = [double]"0.368415"
$arr[27,8] = [double]"0.03446906976034867"
$arr[27,9] = [double]"0.03446906976034867"
$arr[27,10] = 3
$arr[27,11] = 1
$arr[27,12] = [double]"0.8044289999999998"
$arr[27,13] = [double]"2.413287"
$arr[27,14] = [double]"0.005030458199167516"
$arr[27,15] = [double]"0.005030458199167516"
$arr[27,16] = [double]"0.09878790334499998"
$arr[27,17] = [double]"0.8890911301049998"
$arr[27,18] = [double]"0.0001733952145936231"
$arr[27,19] = [double]"0.0001733952145936231"

# row 30: Neutrophils -> Neutrophils
$arr[28,0] = "Neutrophils"
$arr[28,1] = "Spon2"
$arr[28,2] = "Itga4"
$arr[28,3] = "Neutrophils"
$arr[28,4] = 2
$arr[28,5] = [double]"0.6666666666666666"
$arr[28,6] = [double]"0.122805"
$arr[28,7] = [double]"0.368415"
$arr[28,8] = [double]"0.03446906976034867"
$arr[28,9] = [double]"0.03446906976034867"
$arr[28,10] = 3
$arr[28,11] = 1
$arr[28,12] = [double]"78.08909333333334"
$arr[28,13] = [double]"234.26728"
$arr[28,14] = [double]"0.4883264027331488"
$arr[28,15] = [double]"0.4883264027331488"
$arr[28,16] = [double]"9.5897311068"
$arr[28,17] = [double]"86.30757996119999"
$arr[28,18] = [double]"0.01683215684162903"
$arr[28,19] = [double]"0.01683215684162902"

# row 31: Neutrophils -> Resolving-Mac
$arr[29,0] = "Neutrophils"
$arr[29,1] = "Spon2"
$arr[29,2] = "Itga4"
$arr[29,3] = "Resolving-Mac"
$arr[29,4] = 2
$arr[29,5] = [double]"0.6666666666666666"
$arr[29,6] = [double]"0.122805"
$arr[29,7] = [double]"0.368415"
$arr[29,8] = [double]"0.03446906976034867"
$arr[29,9] = [double]"0.03446906976034867"
$arr[29,10] = 3
$arr[29,11] = 1
$arr[29,12] = [double]"23.323433"
$arr[29,13] = [double]"69.970299"
$arr[29,14] = [double]"0.1458519704878668"
$arr[29,15] = [double]"0.1458519704878668"
$arr[29,16] = [double]"2.864234189565"
$arr[29,17] = [double]"25.778107706085"
$arr[29,18] = [double]"0.005027381745430597"
$arr[29,19] = [double]"0.005027381745430597"

$ws.Range("A2:T31").Value = $arr
